$d = $word.ActiveDocument

$d.Content.Find.Execute("11×81=891", $true, $false, $false, $false, $false, $true, 1, $false, "52×20=1040", 2) | Out-Null
$d.Content.Find.Execute("70×31=2170", $true, $false, $false, $false, $false, $true, 1, $false, "51×78=3978", 2) | Out-Null
$d.Content.Find.Execute("25×50=1250", $true, $false, $false, $false, $false, $true, 1, $false, "47×85=3995", 2) | Out-Null
$d.Content.Find.Execute("65×24=1560", $true, $false, $false, $false, $false, $true, 1, $false, "27×11=297", 2) | Out-Null
$d.Content.Find.Execute("49×94=4606", $true, $false, $false, $false, $false, $true, 1, $false, "18×96=1728", 2) | Out-Null
$d.Content.Find.Execute("45×85=3825", $true, $false, $false, $false, $false, $true, 1, $false, "51×27=1377", 2) | Out-Null
$d.Content.Find.Execute("67×11=737", $true, $false, $false, $false, $false, $true, 1, $false, "70×76=5320", 2) | Out-Null
$d.Content.Find.Execute("29×30=870", $true, $false, $false, $false, $false, $true, 1, $false, "41×42=1722", 2) | Out-Null
$d.Content.Find.Execute("48×12=576", $true, $false, $false, $false, $false, $true, 1, $false, "61×41=2501", 2) | Out-Null
$d.Content.Find.Execute("60×24=1440", $true, $false, $false, $false, $false, $true, 1, $false, "53×78=4134", 2) | Out-Null
$d.Content.Find.Execute("44×72=3168", $true, $false, $false, $false, $false, $true, 1, $false, "90×78=7020", 2) | Out-Null
$d.Content.Find.Execute("54×32=1728", $true, $false, $false, $false, $false, $true, 1, $false, "30×33=990", 2) | Out-Null
$d.Content.Find.Execute("68×40=2720", $true, $false, $false, $false, $false, $true, 1, $false, "61×77=4697", 2) | Out-Null
$d.Content.Find.Execute("26×50=1300", $true, $false, $false, $false, $false, $true, 1, $false, "22×41=902", 2) | Out-Null
$d.Content.Find.Execute("53×15=795", $true, $false, $false, $false, $false, $true, 1, $false, "34×47=1598", 2) | Out-Null
$d.Content.Find.Execute("65×36=2340", $true, $false, $false, $false, $false, $true, 1, $false, "49×87=4263", 2) | Out-Null
$d.Content.Find.Execute("80×81=6480", $true, $false, $false, $false, $false, $true, 1, $false, "17×66=1122", 2) | Out-Null
$d.Content.Find.Execute("58×67=3886", $true, $false, $false, $false, $false, $true, 1, $false, "58×17=986", 2) | Out-Null
$d.Content.Find.Execute("24×20=480", $true, $false, $false, $false, $false, $true, 1, $false, "71×15=1065", 2) | Out-Null
$d.Content.Find.Execute("15×96=1440", $true, $false, $false, $false, $false, $true, 1, $false, "39×32=1248", 2) | Out-Null
$d.Content.Find.Execute("56×68=3808", $true, $false, $false, $false, $false, $true, 1, $false, "29×16=464", 2) | Out-Null
$d.Content.Find.Execute("15×97=1455", $true, $false, $false, $false, $false, $true, 1, $false, "98×75=7350", 2) | Out-Null
$d.Content.Find.Execute("12×83=996", $true, $false, $false, $false, $false, $true, 1, $false, "83×20=1660", 2) | Out-Null
$d.Content.Find.Execute("32×90=2880", $true, $false, $false, $false, $false, $true, 1, $false, "77×19=1463", 2) | Out-Null
$d.Content.Find.Execute("62×61=3782", $true, $false, $false, $false, $false, $true, 1, $false, "19×47=893", 2) | Out-Null
